# "Meget lille ændring i forgrunden" (Danish: "Very small change in the
# foreground") - nudges the seven foreground "tool" pictures on slide 2
# (wrench/drill/screwdriver/hammer/pin strips) down-and-right by a tiny
# amount (~22715 EMU right, ~51044 EMU down -> ~1.79pt right, ~4.02pt down).
#
# The shapes are identified by their stable PowerPoint shape Id (not by
# collection index, since that can shift) and are plain siblings in the
# slide's shape tree (not grouped).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Target Left/Top, in points, for each affected shape (by shape Id).
# These correspond exactly to the new EMU offsets:
#   7    -> (1116008, 3150312)
#   23   -> (  22715, 2399527)
#   5    -> (-1950434, 3150313)
#   3    -> ( 891222, 3725194)
#   2062 -> (1644681, 3829222)
#   2064 -> (2598801, 3946916)
#   2058 -> (1998354, 3764781)
$targets = @{
    7    = @{ Left = 87.87464904785156;  Top = 248.0560760498047 }
    23   = @{ Left = 1.788582682609558;  Top = 188.93914794921875 }
    5    = @{ Left = -153.57748413085938; Top = 248.05615234375 }
    3    = @{ Left = 70.17496490478516;  Top = 293.3223876953125 }
    2062 = @{ Left = 129.50244140625;    Top = 301.5135498046875 }
    2064 = @{ Left = 204.6300048828125;  Top = 310.7807922363281 }
    2058 = @{ Left = 157.35072326660156; Top = 296.439453125 }
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    $key = [string]$shape.Id
    if ($targets.ContainsKey($key)) {
        $t = $targets[$key]
        $shape.Left = $t.Left
        $shape.Top = $t.Top
    }
}
